# "Post processing Effects Done" - mark the "Post Processing Effects" WIP
# sub-section (rows 35-37: Post Processing Effects / Glow / Torchlight) as
# Done: fill in the actual "Time Taken" hours, flip the Status label from
# "WIP" to "Done", and re-colour the block from the WIP yellow to the Done
# green (matching the styling already used by the other "Done" sections,
# e.g. row 33/34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Time Taken (hrs) values -------------------------------------------------
# D35 is a SUM formula over D36:D37, so it recalculates automatically.
$ws.Range("D36").Value = 4
$ws.Range("D37").Value = 2

# --- Status: WIP -> Done -----------------------------------------------------
$ws.Range("E35").Value = "Done"

# --- Re-colour the block from "WIP" yellow to "Done" green ------------------
$doneGreen = 5296274   # RGB(146,208,80) / 92D050 - same fill used by other Done rows

$ws.Range("A35").Interior.Color = $doneGreen
$ws.Range("C35:E35").Interior.Color = $doneGreen
$ws.Range("B36:E36").Interior.Color = $doneGreen
$ws.Range("B37:E37").Interior.Color = $doneGreen

# --- View state: scroll down a bit and leave the selection on D42 -----------
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("D42").Select()
